$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 757, shifting existing rows (757-798) down to (758-799)
$ws.Rows(757).Insert()

# Populate the newly inserted row 757 with the new daily-push data
# (force column A to plain text so the date-like string is not
# auto-converted into a date serial number by Excel)
$ws.Range("A757").NumberFormat = "@"
$ws.Range("A757").Value = "2026/02/05"
$ws.Range("B757").Value = "木"
$ws.Range("C757").Value = 5
$ws.Range("D757").Value = 37
